# Applies the "Add files via upload" edit to the kbase workbook:
#  1. Moves the "3.0.3 software issue" FAQ row (old row 47) up to become
#     row 10, updating its date to 45474 and its category tag from
#     "EM" to "EM, Android".
#  2. Fixes up the sheet view (drop the stale topLeftCell / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Move old row 47 ("3.0.3 소프트웨어 이슈" / EM30001) to row 10 ---------

# Remember the source row's height before we start shuffling rows around.
$srcHeight = $ws.Rows.Item(47).RowHeight

# Insert a fresh blank row at position 10; this pushes the current rows
# 10..117 down to 11..118, so the source row we want is now row 48.
$ws.Rows.Item(10).Insert()

# Copy the (now shifted) source row's A:F cells - values AND formatting -
# into the newly inserted row 10.
$ws.Range("A48:F48").Copy($ws.Range("A10:F10"))

# Re-apply the row height, since Range.Copy does not carry row height.
$ws.Rows.Item(10).RowHeight = $srcHeight

# Remove the now-duplicate old row (still sitting at row 48).
$ws.Rows.Item(48).Delete()

# Update the moved row's date and category to their new values.
$ws.Range("B10").Value = 45474
$ws.Range("E10").Value = "EM, Android"

# --- Sheet view tidy-up: drop stale topLeftCell, reset selection to B2 ----

$ws.Activate()
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 1
